$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 66 ("شاش 10 سم"), shifting rows 66-73 down to 67-74.
$ws.Rows.Item(66).Insert()

# Fill in the new row's data (matches the row-66 style pattern: A=6, B-G=7, H-K=8, L-M=9, N=10)
$ws.Range("A66").Value = 63
$ws.Range("B66").Value = "سكاته حصيره "
$ws.Range("H66").Value = "8:0"
$ws.Range("L66").Value = 15
$ws.Range("N66").Value = 1

# Re-create the merges for the new row, matching the other item rows
$ws.Range("B66:G66").Merge()
$ws.Range("H66:K66").Merge()
$ws.Range("L66:M66").Merge()

# Renumber the serial ("م") column for the rows pushed down by the insert
for ($r = 67; $r -le 72; $r++) {
    $ws.Range("A$r").Value = $r - 3
}

# Update the running total (K column, now shifted to row 73)
$ws.Range("K73").Value = 5563.17
